$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing ticker lists (columns B, C, E) for rows 2-17 ---

# Row 2
$ws.Range("B2").Value = "NSE:BRITANNIA"
$ws.Range("C2").Value = "NSE:ADANIENT"
$ws.Range("E2").Value = "NSE:APOLLOTYRE"

# Row 3
$ws.Range("B3").Value = "NSE:CUMMINSIND"
$ws.Range("C3").Value = "NSE:ARMANFIN"
$ws.Range("E3").ClearContents()
$ws.Range("E3").Font.Bold = $false

# Row 4
$ws.Range("B4").Value = "NSE:DMART"
$ws.Range("C4").Value = "NSE:AROGRANITE"
$ws.Range("E4").ClearContents()
$ws.Range("E4").Font.Bold = $false

# Row 5
$ws.Range("B5").Value = "NSE:MONQ50"
$ws.Range("C5").Value = "NSE:ATFL"
$ws.Range("E5").ClearContents()
$ws.Range("E5").Font.Bold = $false

# Row 6
$ws.Range("C6").Value = "NSE:CGCL"
$ws.Range("E6").ClearContents()
$ws.Range("E6").Font.Bold = $false

# Row 7
$ws.Range("C7").Value = "NSE:DIAMINESQ"
$ws.Range("E7").ClearContents()
$ws.Range("E7").Font.Bold = $false

# Row 8-17 column C only
$ws.Range("C8").Value = "NSE:DODLA"
$ws.Range("C9").Value = "NSE:DYNPRO"
$ws.Range("C10").Value = "NSE:FDC"
$ws.Range("C11").Value = "NSE:GRWRHITECH"
$ws.Range("C12").Value = "NSE:GTPL"
$ws.Range("C13").Value = "NSE:HONASA"
$ws.Range("C14").Value = "NSE:ISGEC"
$ws.Range("C15").Value = "NSE:MANYAVAR"
$ws.Range("C16").Value = "NSE:MGL"
$ws.Range("C17").Value = "NSE:MINDACORP"

# --- Append new rows 18-22, following the same layout as the existing rows:
#     column A holds a styled running index, column C the ticker, and
#     columns B/D/E/F exist as (empty) cells. ---
$newRows = @(
    @{ Row = 18; Num = 16; Ticker = "NSE:NCLIND" },
    @{ Row = 19; Num = 17; Ticker = "NSE:NUVOCO" },
    @{ Row = 20; Num = 18; Ticker = "NSE:ONMOBILE" },
    @{ Row = 21; Num = 19; Ticker = "NSE:POKARNA" },
    @{ Row = 22; Num = 20; Ticker = "NSE:SAKSOFT" }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the style (border/bold/centered alignment) used by column A in the
    # existing data rows, then set the running index value.
    $ws.Range("A17").Copy($ws.Range("A$rowNum"))
    $ws.Range("A$rowNum").Value = $r.Num

    # Materialize the other blank cells of the row (B, D, E, F) without
    # introducing any new cell style.
    $ws.Range("B$rowNum").Font.Bold = $false
    $ws.Range("C$rowNum").Value = $r.Ticker
    $ws.Range("D$rowNum").Font.Bold = $false
    $ws.Range("E$rowNum").Font.Bold = $false
    $ws.Range("F$rowNum").Font.Bold = $false
}
